$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains exact text formatting (e.g. trailing zeros)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '56.674.29'
$ws.Range('E2').Value = '  +2.78%  '
$ws.Range('D3').Value = '3.004.33'
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '510.23'
$ws.Range('E5').Value = '  +6.75%  '
$ws.Range('D6').Value = '139.50'
$ws.Range('E6').Value = '  +7.80%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.433'
$ws.Range('E8').Value = '  +5.36%  '
$ws.Range('D9').Value = '7.56'
$ws.Range('E9').Value = '  +11.81%  '
$ws.Range('E10').Value = '  +10.14%  '
$ws.Range('D11').Value = '0.355'
$ws.Range('E11').Value = '  +4.24%  '
$ws.Range('E12').Value = '  +3.80%  '
$ws.Range('D13').Value = '3.518.84'
$ws.Range('E13').Value = '  +2.65%  '
$ws.Range('D14').Value = '25.65'
$ws.Range('E14').Value = '  +7.81%  '
$ws.Range('E15').Value = '  +13.84%  '
$ws.Range('D16').Value = '56.758.37'
$ws.Range('E16').Value = '  +3.17%  '
$ws.Range('D17').Value = '3.001.81'
$ws.Range('E17').Value = '  +2.89%  '
$ws.Range('D18').Value = '5.95'
$ws.Range('E18').Value = '  +8.40%  '
$ws.Range('D19').Value = '12.50'
$ws.Range('E19').Value = '  +6.65%  '
$ws.Range('D20').Value = '7.85'
$ws.Range('E20').Value = '  +8.18%  '
$ws.Range('D21').Value = '329.05'
$ws.Range('E21').Value = '  +7.15%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = '0.483'
$ws.Range('E23').Value = '  +7.11%  '
$ws.Range('D24').Value = '62.88'
$ws.Range('E24').Value = '  +5.61%  '
$ws.Range('E25').Value = '  +13.05%  '
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('D27').Value = '0.0₃0913'
$ws.Range('E27').Value = '  +10.71%  '
$ws.Range('E28').Value = '  +5.36%  '
$ws.Range('D29').Value = '7.08'
$ws.Range('E29').Value = '  +11.09%  '
$ws.Range('D30').Value = '1.26'
$ws.Range('E30').Value = '  +10.03%  '
$ws.Range('E31').Value = '  +8.35%  '
$ws.Range('D32').Value = '20.64'
$ws.Range('E32').Value = '  +8.32%  '
$ws.Range('D33').Value = '155.28'
$ws.Range('E33').Value = '  +6.73%  '
$ws.Range('D34').Value = '4.58'
$ws.Range('E34').Value = '  +7.39%  '
$ws.Range('D35').Value = '5.67'
$ws.Range('E35').Value = '  +3.26%  '
$ws.Range('D36').Value = '1.27'
$ws.Range('E36').Value = '  +2.85%  '
$ws.Range('D37').Value = '0.0681'
$ws.Range('E37').Value = '  +8.19%  '
$ws.Range('D38').Value = '24.29'
$ws.Range('E38').Value = '  +3.29%  '
$ws.Range('D39').Value = '3.037.72'
$ws.Range('E39').Value = '  +2.83%  '
$ws.Range('D40').Value = '36.95'
$ws.Range('E40').Value = '  +3.31%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('E42').Value = '  +4.71%  '
$ws.Range('D43').Value = '2.267.41'
$ws.Range('E43').Value = '  +9.35%  '
$ws.Range('E44').Value = '  +2.26%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '3.66'
$ws.Range('E45').Value = '  +5.61%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').Value = '1.41'
$ws.Range('E46').Value = '  +4.52%  '
$ws.Range('D47').Value = '1.97'
$ws.Range('E47').Value = '  +20.43%  '
$ws.Range('E48').Value = '  +8.21%  '
$ws.Range('E49').Value = '  +6.68%  '
$ws.Range('D50').Value = '19.38'
$ws.Range('E50').Value = '  +6.04%  '
$ws.Range('E51').Value = '  +8.42%  '
